# "Updated Columns in DT Sheet"
# The DT sheet (5th tab) was empty; add the exam-summary header row
# (Year / Department / ExamTitle / NumSub) and size columns B & C to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DT")

$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Department"
$ws.Range("C1").Value = "ExamTitle"
$ws.Range("D1").Value = "NumSub"

# Custom widths for the Department / ExamTitle columns.
$ws.Columns.Item(2).ColumnWidth = 11.26
$ws.Columns.Item(3).ColumnWidth = 12.59

# DT becomes the active/selected sheet, with D1 the selected cell.
$ws.Activate()
$ws.Range("D1").Select()
